# Updated symbol list with GitHub Actions - apply cell value changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '308.66'
    'E2' = '-0.69%'
    'D3' = '37.72'
    'E3' = '0.43%'
    'D4' = '5.133'
    'E4' = '0.88%'
    'D5' = '0.07855'
    'E5' = '1.10%'
    'B6' = 'GateToken'
    'C6' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D6' = '4.429'
    'E6' = '1.93%'
    'B7' = 'FTXToken'
    'C7' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D7' = '1.902'
    'E7' = '0.24%'
    'B8' = 'KuCoinToken'
    'C8' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D8' = '8.230'
    'E8' = '0.29%'
    'B9' = 'BTSEToken'
    'C9' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D9' = '2.999'
    'E9' = '2.27%'
    'B10' = 'MXToken'
    'C10' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D10' = '0.9350'
    'E10' = '1.93%'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D11' = '0.1090'
    'E11' = '-9.31%'
    'B12' = 'WazirX'
    'C12' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D12' = '0.1947'
    'E12' = '1.13%'
    'B13' = 'MandalaExchangeToken'
    'C13' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D13' = '0.09007'
    'E13' = '1.02%'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '0.03340'
    'E14' = '-2.17%'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D15' = '0.09597'
    'E15' = '-1.06%'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '0.001384'
    'E16' = '1.20%'
    'B17' = 'TigerCash'
    'C17' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D17' = '0.005714'
    'E17' = '-1.70%'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D18' = '3.614'
    'E18' = '1.72%'
    'D19' = '0.3412'
    'E19' = '0.57%'
    'D20' = '6.228'
    'E20' = '23.76%'
    'D21' = '0.1280'
    'E21' = '0.16%'
    'E22' = '-10.53%'
    'D23' = '0.04400'
    'E23' = '0.69%'
    'E24' = '1.82%'
    'D25' = '0.004566'
    'E25' = '7.31%'
    'D39' = '0.02221'
    'E39' = '5.13%'
    'D40' = '0.05053'
    'E40' = '2.19%'
    'D41' = '0.007467'
    'E41' = '-2.25%'
    'D42' = '0.1351'
    'E42' = '0.64%'
    'D43' = '0.008736'
    'E43' = '-11.52%'
    'D44' = '0.002112'
    'E44' = '2.62%'
    'D45' = '0.008054'
    'E45' = '-16.01%'
    'E46' = '-2.04%'
    'E47' = '0.13%'
    'E48' = '-5.95%'
    'E49' = '-40.73%'
    'E50' = '0.13%'
    'E51' = '0.13%'
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.Style = $origStyle
}
